$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.869.33"
$ws.Range("E2").Value = "  -2.44%  "

# Row 3
$ws.Range("D3").Value = "3.228.40"
$ws.Range("E3").Value = "  -5.40%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'580.84"
$ws.Range("E5").Value = "  -4.86%  "

# Row 6
$ws.Range("D6").Value = "'142.53"
$ws.Range("E6").Value = "  -14.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "3.222.28"
$ws.Range("E8").Value = "  -5.48%  "

# Row 9
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -11.21%  "

# Row 10
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  -14.92%  "

# Row 11
$ws.Range("D11").Value = "'6.45"
$ws.Range("E11").Value = "  -6.07%  "

# Row 12
$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  -12.73%  "

# Row 13
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -11.50%  "

# Row 14
$ws.Range("D14").Value = "'36.51"
$ws.Range("E14").Value = "  -16.27%  "

# Row 15
$ws.Range("D15").Value = "3.744.20"
$ws.Range("E15").Value = "  -5.59%  "

# Row 16
$ws.Range("D16").Value = "67.037.76"
$ws.Range("E16").Value = "  -2.26%  "

# Row 17
$ws.Range("D17").Value = "3.235.83"
$ws.Range("E17").Value = "  -5.36%  "

# Row 18
$ws.Range("E18").Value = "  -6.51%  "

# Row 19
$ws.Range("D19").Value = "'6.83"
$ws.Range("E19").Value = "  -15.05%  "

# Row 20
$ws.Range("D20").Value = "'504.11"
$ws.Range("E20").Value = "  -12.29%  "

# Row 21
$ws.Range("D21").Value = "'14.43"
$ws.Range("E21").Value = "  -14.77%  "

# Row 22
$ws.Range("D22").Value = "'0.726"
$ws.Range("E22").Value = "  -13.44%  "

# Row 23
$ws.Range("D23").Value = "'7.41"
$ws.Range("E23").Value = "  -16.67%  "

# Row 24
$ws.Range("D24").Value = "'82.52"
$ws.Range("E24").Value = "  -12.45%  "

# Row 25
$ws.Range("D25").Value = "'12.88"
$ws.Range("E25").Value = "  -13.69%  "

# Row 26
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("D27").Value = "'3.11"
$ws.Range("E27").Value = "  -13.55%  "

# Row 28
$ws.Range("D28").Value = "'2.06"
$ws.Range("E28").Value = "  -13.23%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'27.94"
$ws.Range("E29").Value = "  -13.41%  "

# Row 30
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.66"
$ws.Range("E30").Value = "  -9.34%  "

# Row 31
$ws.Range("E31").Value = "  -5.57%  "

# Row 32
$ws.Range("D32").Value = "'2.52"
$ws.Range("E32").Value = "  -7.82%  "

# Row 33
$ws.Range("E33").Value = "  +0.12%  "

# Row 34
$ws.Range("D34").Value = "'6.19"
$ws.Range("E34").Value = "  -19.92%  "

# Row 35
$ws.Range("D35").Value = "'5.42"
$ws.Range("E35").Value = "  -16.34%  "

# Row 36
$ws.Range("D36").Value = "'54.38"
$ws.Range("E36").Value = "  -2.95%  "

# Row 37
$ws.Range("D37").Value = "'491.08"
$ws.Range("E37").Value = "  -15.62%  "

# Row 38
$ws.Range("D38").Value = "'0.0421"
$ws.Range("E38").Value = "  -8.60%  "

# Row 39
$ws.Range("D39").Value = "'0.0818"
$ws.Range("E39").Value = "  -13.22%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.122"
$ws.Range("E40").Value = "  -12.70%  "

# Row 41
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.59"
$ws.Range("E41").Value = "  -17.10%  "

# Row 42
$ws.Range("D42").Value = "2.856.60"
$ws.Range("E42").Value = "  -10.66%  "

# Row 43
$ws.Range("D43").Value = "'2.57"
$ws.Range("E43").Value = "  -14.68%  "

# Row 44
$ws.Range("D44").Value = "'0.254"
$ws.Range("E44").Value = "  -12.47%  "

# Row 45
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").Value = "'2.09"
$ws.Range("E46").Value = "  -10.76%  "

# Row 47
$ws.Range("D47").Value = "'25.36"
$ws.Range("E47").Value = "  -18.37%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'121.60"
$ws.Range("E48").Value = "  -7.53%  "

# Row 49
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0532"
$ws.Range("E49").Value = "  -19.97%  "

# Row 50
$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "  -12.19%  "

# Row 51
$ws.Range("D51").Value = "'2.16"
$ws.Range("E51").Value = "  -21.45%  "
